$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.744.45"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "'3.706.16"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'677.95"
$ws.Range("E5").Value = "  -1.15%  "
$ws.Range("D6").Value = "'162.61"
$ws.Range("E6").Value = "  +1.41%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.149"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").Value = "'7.16"
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("D11").Value = "'0.444"
$ws.Range("E11").Value = "  +1.88%  "
$ws.Range("D12").Value = "'0.0000236"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").Value = "'32.92"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("D14").Value = "'3.706.60"
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "'69.738.81"
$ws.Range("E15").Value = "  +0.42%  "
$ws.Range("E16").Value = "  +1.82%  "
$ws.Range("D17").Value = "'16.13"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "'6.52"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "'473.94"
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").Value = "'9.86"
$ws.Range("E20").Value = "  -1.40%  "
$ws.Range("D21").Value = "'0.656"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").Value = "'80.46"
$ws.Range("D23").Value = "'3.853.34"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").Value = "'0.0000128"
$ws.Range("E24").Value = "  +3.08%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").Value = "'9.16"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("D31").Value = "'6.63"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("E32").Value = "  +5.02%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'27.03"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "'3.693.78"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").Value = "'8.57"
$ws.Range("E36").Value = "  +4.44%  "
$ws.Range("D37").Value = "'6.23"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("E41").Value = "  +0.56%  "
$ws.Range("D42").Value = "'0.946"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'167.41"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "'47.01"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("D45").Value = "'2.79"
$ws.Range("E45").Value = "  +1.91%  "
$ws.Range("D46").Value = "'0.000282"
$ws.Range("E46").Value = "  -0.80%  "
$ws.Range("D47").Value = "'28.20"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("D50").Value = "'7.94"
$ws.Range("E50").Value = "  +2.02%  "
$ws.Range("D51").Value = "'0.268"
$ws.Range("E51").Value = "  +2.06%  "
